$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "fenerbahçe en son ne zaman şampiyonlar ligi'ne katıldı."
